$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.434.67'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.51%  '
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.238.58'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.51%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.629'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +0.69%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '74.03'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -4.14%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.621'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '43.38'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +3.96%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0960'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.13'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.09%  '
$ws.Range('E13').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '14.46'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -2.70%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.855'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.266.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +0.16%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '42.273.03'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.35%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +11.50%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.64%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.02'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +39.45%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '231.54'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -4.97%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.79'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +3.92%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.73%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -0.31%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.23'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.01%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '166.73'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.64%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '21.00'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.90%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.87'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +18.52%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0811'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -2.44%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.36%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '30.15'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -8.75%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.49%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.45%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +1.69%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '13.20'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -7.15%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.90%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.47%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '63.50'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +3.66%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.92'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +2.73%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '105.85'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -6.94%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +2.66%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.996'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +0.87%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.38'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.71%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.78%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.74'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.38%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -2.60%  '
$ws.Range('E51').Style = 'Normal'
